$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 118
$ws.Range("I2").Value = 118
$ws.Range("K2").Value = 118
$ws.Range("M2").Value = -5

$ws.Range("H12").Value = 3769.1428
$ws.Range("I12").Value = 4044.6538
$ws.Range("J12").Value = 187.5
$ws.Range("K12").Value = 4044.6538
$ws.Range("L12").Value = 187.5
$ws.Range("M12").Value = -3874.6538
$ws.Range("N12").Value = -527.5

$ws.Range("H40").Value = 3010.8125
$ws.Range("I40").Value = 2198
$ws.Range("K40").Value = 2198
$ws.Range("M40").Value = -2023

$ws.Range("H70").Value = 4473
$ws.Range("J70").Value = 4773.8
$ws.Range("L70").Value = 14321.4
$ws.Range("N70").Value = -14861.4

$ws.Range("H73").Value = 4473
$ws.Range("J73").Value = 4773.8
$ws.Range("L73").Value = 14321.4
$ws.Range("N73").Value = -16193.4

$ws.Range("H74").Value = 3966.6667
$ws.Range("I74").Value = 3975
$ws.Range("J74").Value = 3950
$ws.Range("K74").Value = 3975
$ws.Range("L74").Value = 3950
$ws.Range("M74").Value = -3039
$ws.Range("N74").Value = -5822

$ws.Range("H77").Value = 3966.6667
$ws.Range("I77").Value = 3975
$ws.Range("J77").Value = 3950
$ws.Range("K77").Value = 19875
$ws.Range("L77").Value = 19750
$ws.Range("M77").Value = -15195
$ws.Range("N77").Value = -29110

$ws.Range("H116").Value = 4511.25
$ws.Range("I116").Value = 4342.5293
$ws.Range("J116").Value = 5467.3335
$ws.Range("K116").Value = 4342.5293
$ws.Range("L116").Value = 5467.3335
$ws.Range("M116").Value = -900.5293000000001
$ws.Range("N116").Value = -12351.3335

$ws.Range("H119").Value = 995
$ws.Range("I119").Value = 995
$ws.Range("K119").Value = 2985
$ws.Range("M119").Value = 1853

$ws.Range("H129").Value = 2059
$ws.Range("I129").Value = 2385
$ws.Range("J129").Value = 1896
$ws.Range("K129").Value = 7155
$ws.Range("L129").Value = 5688
$ws.Range("M129").Value = -2155
$ws.Range("N129").Value = -15688

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

$ws.Range("H137").Value = 1328.5625
$ws.Range("I137").Value = 1034.381
$ws.Range("K137").Value = 3103.143
$ws.Range("M137").Value = -553.143

$ws.Range("H138").Value = 3871.9666
$ws.Range("J138").Value = 3924.3044
$ws.Range("L138").Value = 11772.9132
$ws.Range("N138").Value = -22052.9132

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H45").Value = 5246.154
$ws.Range("I45").Value = 8714.833000000001
$ws.Range("J45").Value = 2273
$ws.Range("K45").Value = 8714.833000000001
$ws.Range("L45").Value = 2273
$ws.Range("M45").Value = -8337.833000000001
$ws.Range("N45").Value = -3027

$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H132").Value = 3335814
$ws.Range("I132").Value = 3573879.2
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 10721637.6
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -10719107.6
$ws.Range("N132").Value = -13760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 134771
$ws.Range("I107").Value = 797
$ws.Range("J107").Value = 251998.25
$ws.Range("K107").Value = 797
$ws.Range("L107").Value = 251998.25
$ws.Range("M107").Value = 1123
$ws.Range("N107").Value = -255838.25

$ws.Range("H134").Value = 26318206
$ws.Range("I134").Value = 27780190
$ws.Range("J134").Value = 2499
$ws.Range("K134").Value = 83340570
$ws.Range("L134").Value = 7497
$ws.Range("M134").Value = -83338035
$ws.Range("N134").Value = -12567

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 184.375
$ws.Range("I7").Value = 237.8
$ws.Range("J7").Value = 95.333336
$ws.Range("K7").Value = 237.8
$ws.Range("L7").Value = 95.333336
$ws.Range("M7").Value = -124.8
$ws.Range("N7").Value = -321.333336

$ws.Range("H22").Value = 540.875
$ws.Range("I22").Value = 540.875
$ws.Range("K22").Value = 540.875
$ws.Range("M22").Value = -190.875

$ws.Range("H62").Value = 7626.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7626.5
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = 7626.5
$ws.Range("N62").Value = -8874.5

$ws.Range("H65").Value = 7626.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7626.5
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = 38132.5
$ws.Range("N65").Value = -44372.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2108.5
$ws.Range("I68").Value = 1241.6666
$ws.Range("J68").Value = 2261.4707
$ws.Range("K68").Value = 3724.9998
$ws.Range("L68").Value = 6784.4121
$ws.Range("M68").Value = -2913.9998
$ws.Range("N68").Value = -8406.4121

$ws.Range("H71").Value = 2108.5
$ws.Range("I71").Value = 1241.6666
$ws.Range("J71").Value = 2261.4707
$ws.Range("K71").Value = 11174.9994
$ws.Range("L71").Value = 20353.2363
$ws.Range("M71").Value = -7118.999400000001
$ws.Range("N71").Value = -28465.2363

$ws.Range("H118").Value = 1120
$ws.Range("I118").Value = 1120
$ws.Range("K118").Value = 3360
$ws.Range("M118").Value = -2117

$ws.Range("H122").Value = 928.5
$ws.Range("I122").Value = 900.5
$ws.Range("K122").Value = 8104.5
$ws.Range("M122").Value = -5654.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3491.2942
$ws.Range("I61").Value = 3491.2942
$ws.Range("K61").Value = 3491.2942
$ws.Range("M61").Value = -3289.2942

$ws.Range("H93").Value = 2209.5386
$ws.Range("I93").Value = 1958.2222
$ws.Range("K93").Value = 1958.2222
$ws.Range("M93").Value = -710.2221999999999

$ws.Range("H113").Value = 3491.2942
$ws.Range("I113").Value = 3491.2942
$ws.Range("K113").Value = 3491.2942
$ws.Range("M113").Value = -1321.2942

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1893.1428
$ws.Range("I122").Value = 1883.6666
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 5650.9998
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -3200.9998
$ws.Range("N122").Value = -10750

$ws.Range("H126").Value = 2133.2778
$ws.Range("I126").Value = 2019.6923
$ws.Range("K126").Value = 6059.0769
$ws.Range("M126").Value = -3589.0769

$ws.Range("H132").Value = 33334706
$ws.Range("I132").Value = 38462884
$ws.Range("J132").Value = 1542
$ws.Range("K132").Value = 115388652
$ws.Range("L132").Value = 4626
$ws.Range("M132").Value = -115386122
$ws.Range("N132").Value = -9686
